# Energia_Primaria_Lisboa.xlsx - data refresh ("Add files via upload")
#
# The underlying G-column totals were recomputed upstream and a handful of
# rows changed value; the active sheet's selection was also left on a
# different range when the file was re-saved. Reproduce both here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Energia_Primaria")

# --- Updated totals in column G (rows scattered through the table) ---
$ws.Range("G9").Value  = 378091.86347555893
$ws.Range("G19").Value = 354171.50838907313
$ws.Range("G24").Value = 350683.44765434955
$ws.Range("G29").Value = 361684.42839035933
$ws.Range("G34").Value = 344586.54493832408
$ws.Range("G36").Value = 1918.3988340409351
$ws.Range("G37").Value = 16942.006527913905
$ws.Range("G38").Value = 22442.227475102762
$ws.Range("G39").Value = 371858.81673632929
$ws.Range("G40").Value = 120572.562235836
$ws.Range("G44").Value = 376095.50883206091
$ws.Range("G49").Value = 366767.9169788835
$ws.Range("G54").Value = 385884.80354058166

# --- Selection left active on C6:H55 (anchor C6) when the author saved ---
$null = $ws.Range("C6:H55").Select()
